$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct existing values for rows 174-176 (columns B and D)
$ws.Range("B174").Value = 13711
$ws.Range("D174").Value = 1349

$ws.Range("B175").Value = 12516
$ws.Range("D175").Value = 1300

$ws.Range("B176").Value = 12984
$ws.Range("D176").Value = 1374

# Add the new row 177 with the August 2021 data
$ws.Range("A177").NumberFormat = "@"
$ws.Range("A177").Value = "01-08-2021"
$ws.Range("A177").Style = "Normal"
$ws.Range("B177").Value = 12251
$ws.Range("C177").Value = 1348
$ws.Range("D177").Value = 1307
$ws.Range("E177").Value = 144
